$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Sheet1")

# Fix typo: "Savana Pack" -> "Savanna Pack"
$ws.Range("B7").Value = "Savanna Pack"

# Populate Bugs Pack tags (E11) which was left as a placeholder 0
$ws.Range("E11").Value = "animal flashcards, kids flashcards, printable learning, Montessori cards, educational cards, digital download, PDF flashcards, homeschool activities, preschool learning, bugs flashcards, insects learning, creepy crawlies, early learning"

$wb.Save()
